$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 3-9 (individual cell changes per diff) ---

# Row 3
$ws.Range("D3").Value = 44418
$ws.Range("J3").Value = 12

# Row 4
$ws.Range("D4").Value = 44340
$ws.Range("J4").Value = 25

# Row 5
$ws.Range("D5").Value = 44421
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 15000
$ws.Range("P5").Value = 600

# Row 6
$ws.Range("D6").Value = 44432

# Row 7
$ws.Range("D7").Value = 44446
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 13000
$ws.Range("P7").Value = 520

# Row 8
$ws.Range("D8").Value = 44435
$ws.Range("J8").Value = 15
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 14000
$ws.Range("P8").Value = 560

# Row 9
$ws.Range("D9").Value = 44435

# --- Append new row 10 ---
$ws.Range("A10").Value = 12
$ws.Range("B10").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 44376
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 100112026
$ws.Range("G10").Value = "Haba"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 15
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 12000
$ws.Range("N10").Value = "$/saco 25 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 480
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
